# Update rotation values (column D) for parts updated for LCSC orders
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPL")

$ws.Range("D9").Value = 90
$ws.Range("D11").Value = 180
$ws.Range("D12").Value = 90
$ws.Range("D17").Value = 270

# Update the active selection on the sheet to D16
$ws.Activate()
$ws.Range("D16").Select()
